$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.28294533333333
$ws.Range("H2").Value = 57.848836
$ws.Range("I2").Value = 0.04564777115344932
$ws.Range("J2").Value = 0.04564777115344931
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 3241.652610464772
$ws.Range("R2").Value = 29174.87349418295
$ws.Range("S2").Value = 0.01362212687948521
$ws.Range("T2").Value = 0.01362212687948521
$ws.Range("G3").Value = 19.28294533333333
$ws.Range("H3").Value = 57.848836
$ws.Range("I3").Value = 0.04564777115344932
$ws.Range("J3").Value = 0.04564777115344931
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 3143.240331352784
$ws.Range("R3").Value = 28289.16298217505
$ws.Range("S3").Value = 0.01320857715233829
$ws.Range("T3").Value = 0.01320857715233828
$ws.Range("G4").Value = 19.28294533333333
$ws.Range("H4").Value = 57.848836
$ws.Range("I4").Value = 0.04564777115344932
$ws.Range("J4").Value = 0.04564777115344931
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 3200.844164677027
$ws.Range("R4").Value = 28807.59748209324
$ws.Range("S4").Value = 0.01345064094527971
$ws.Range("T4").Value = 0.0134506409452797
$ws.Range("G5").Value = 19.28294533333333
$ws.Range("H5").Value = 57.848836
$ws.Range("I5").Value = 0.04564777115344932
$ws.Range("J5").Value = 0.04564777115344931
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 1277.046497754857
$ws.Range("R5").Value = 11493.41847979371
$ws.Range("S5").Value = 0.005366426176346119
$ws.Range("T5").Value = 0.005366426176346117
$ws.Range("G6").Value = 268.8003336666666
$ws.Range("H6").Value = 806.401001
$ws.Range("I6").Value = 0.6363206400827226
$ws.Range("J6").Value = 0.6363206400827226
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 45187.97767984571
$ws.Range("R6").Value = 406691.7991186114
$ws.Range("S6").Value = 0.1898896764554757
$ws.Range("T6").Value = 0.1898896764554757
$ws.Range("G7").Value = 268.8003336666666
$ws.Range("H7").Value = 806.401001
$ws.Range("I7").Value = 0.6363206400827226
$ws.Range("J7").Value = 0.6363206400827226
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.1841248774207198
$ws.Range("T7").Value = 0.1841248774207198
$ws.Range("G8").Value = 268.8003336666666
$ws.Range("H8").Value = 806.401001
$ws.Range("I8").Value = 0.6363206400827226
$ws.Range("J8").Value = 0.6363206400827226
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 44619.11625050784
$ws.Range("R8").Value = 401572.0462545706
$ws.Range("S8").Value = 0.1874991974318228
$ws.Range("T8").Value = 0.1874991974318228
$ws.Range("G9").Value = 268.8003336666666
$ws.Range("H9").Value = 806.401001
$ws.Range("I9").Value = 0.6363206400827226
$ws.Range("J9").Value = 0.6363206400827226
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 17801.76828645369
$ws.Range("R9").Value = 160215.9145780832
$ws.Range("S9").Value = 0.07480688877470432
$ws.Range("T9").Value = 0.07480688877470432
$ws.Range("G10").Value = 56.43559133333333
$ws.Range("H10").Value = 169.306774
$ws.Range("I10").Value = 0.1335977939863952
$ws.Range("J10").Value = 0.1335977939863952
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 9487.377514501228
$ws.Range("R10").Value = 85386.39763051106
$ws.Range("S10").Value = 0.03986801665264841
$ws.Range("T10").Value = 0.03986801665264842
$ws.Range("G11").Value = 56.43559133333333
$ws.Range("H11").Value = 169.306774
$ws.Range("I11").Value = 0.1335977939863952
$ws.Range("J11").Value = 0.1335977939863952
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 9199.353300869023
$ws.Range("R11").Value = 82794.17970782121
$ws.Range("S11").Value = 0.03865767647930723
$ws.Range("T11").Value = 0.03865767647930723
$ws.Range("G12").Value = 56.43559133333333
$ws.Range("H12").Value = 169.306774
$ws.Range("I12").Value = 0.1335977939863952
$ws.Range("J12").Value = 0.1335977939863952
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 9367.943023057405
$ws.Range("R12").Value = 84311.48720751665
$ws.Range("S12").Value = 0.03936612703283462
$ws.Range("T12").Value = 0.03936612703283462
$ws.Range("G13").Value = 56.43559133333333
$ws.Range("H13").Value = 169.306774
$ws.Range("I13").Value = 0.1335977939863952
$ws.Range("J13").Value = 0.1335977939863952
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 3737.544914177236
$ws.Range("R13").Value = 33637.90422759512
$ws.Range("S13").Value = 0.01570597382160492
$ws.Range("T13").Value = 0.01570597382160492
$ws.Range("G14").Value = 77.91019566666667
$ws.Range("H14").Value = 233.730587
$ws.Range("I14").Value = 0.184433794777433
$ws.Range("J14").Value = 0.1844337947774329
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 13097.46954102955
$ws.Range("R14").Value = 117877.2258692659
$ws.Range("S14").Value = 0.05503840581564262
$ws.Range("T14").Value = 0.05503840581564262
$ws.Range("G15").Value = 77.91019566666667
$ws.Range("H15").Value = 233.730587
$ws.Range("I15").Value = 0.184433794777433
$ws.Range("J15").Value = 0.1844337947774329
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 12699.84771567678
$ws.Range("R15").Value = 114298.629441091
$ws.Range("S15").Value = 0.05336751272317417
$ws.Range("T15").Value = 0.05336751272317417
$ws.Range("G16").Value = 77.91019566666667
$ws.Range("H16").Value = 233.730587
$ws.Range("I16").Value = 0.184433794777433
$ws.Range("J16").Value = 0.1844337947774329
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 12932.5884017007
$ws.Range("R16").Value = 116393.2956153063
$ws.Range("S16").Value = 0.05434553953110587
$ws.Range("T16").Value = 0.05434553953110587
$ws.Range("G17").Value = 77.91019566666667
$ws.Range("H17").Value = 233.730587
$ws.Range("I17").Value = 0.184433794777433
$ws.Range("J17").Value = 0.1844337947774329
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 5159.737830274352
$ws.Range("R17").Value = 46437.64047246917
$ws.Range("S17").Value = 0.02168233670751031
$ws.Range("T17").Value = 0.0216823367075103
